# Add a new "dictation" column (I) to Sheet1, mirroring the existing
# question-type columns (e.g. B = "Listening"/"Reading"), and set it for
# every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I1").Value = "dictation"
$ws.Range("I2").Value = "dictation"
$ws.Range("I3").Value = "dictation"
$ws.Range("I4").Value = "dictation"

$ws.Range("I8").Select()
